$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Monday (row block 4-10)
# ---------------------------------------------------------------------------

# D4/E4 -> cleared out (was "Разр. ПО для моб. Платформ ФМЕ ПОИТ3 ЗФПО" / "210")
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = ""

# D5/E5 -> "Резерв" -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D5").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "364"

# D6/E6 -> "СПП (ЗО) 2 курс колония" -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364 (was 115)
$ws.Range("D6").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "364"

# ---------------------------------------------------------------------------
# Tuesday (row block 14-20)
# ---------------------------------------------------------------------------

# D14/E14 -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D14").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "364"

# D15/E15 -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D15").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "364"

# D16/E16 -> "Разр. ПО для моб. Платформ ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D16").Value = "Разр. ПО для моб. Платформ ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "364"

# ---------------------------------------------------------------------------
# Wednesday (row block 24-30)
# ---------------------------------------------------------------------------

# D24/E24 -> "Разр. ПО для моб. Платформ ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D24").Value = "Разр. ПО для моб. Платформ ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "364"

# ---------------------------------------------------------------------------
# Thursday (row block 34-40)
# ---------------------------------------------------------------------------

# D34/E34 -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D34").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "364"

# D35/E35 -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D35").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "364"

# D36/E36 -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D36").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "364"

# ---------------------------------------------------------------------------
# Friday (row block 44-50)
# ---------------------------------------------------------------------------

# D46/E46 -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D46").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "364"

# D47/E47 -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D47").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "364"

# D48/E48 -> "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО" / 364
$ws.Range("D48").Value = "СПП (ЗО) 2 курс ФМЕ ПОИТ3 ЗФПО"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "364"
